$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table update (GitHub Actions data refresh)
# D column = Price (text-formatted), E column = Volume(1h) % change (text-formatted)

$ws.Range("D2").Value = "24.738.22"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.701.25"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3944"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4041"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08882"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.465"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.189"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.33%  "

$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "1.706.45"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07053"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.071"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.33%  "

$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("E23").Value = "  +3.95%  "

$ws.Range("D24").Value = "24.737.02"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.370"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.686"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.179"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09094"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.667"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.075"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.991"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2758"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02782"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.463"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7708"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7193"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("E45").Value = "  +1.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.222"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.343"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "91.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07987"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
